$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextLooksLikeNumber($cell, $text) {
    # Assigning a numeric-looking string via .Value auto-coerces to a number.
    # Route it through a text-producing formula, then Copy/PasteSpecial
    # (values only) to freeze it back down to a literal text cell without
    # picking up any incidental formatting.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 4 - man1.jpg observation, lat/lon
$ws.Range("T4").Value = "observations/portrait/man1.jpg"
$ws.Range("U4").Value = 37.785834000000001
Set-TextLooksLikeNumber $ws.Range("V4") " -122.406417"

# Row 5 - man2.jpg observation, lat/lon
$ws.Range("T5").Value = "observations/portrait/man2.jpg"
Set-TextLooksLikeNumber $ws.Range("U5") " 37.767087"
Set-TextLooksLikeNumber $ws.Range("V5") " -122.419977"

# Row 6 - woman1.jpg observation, lat/lon (T6 gets the new black-font style)
$ws.Range("T6").Value = "observations/portrait/woman1.jpg"
$ws.Range("T6").Font.Color = 0
Set-TextLooksLikeNumber $ws.Range("U6") " 37.767122"
Set-TextLooksLikeNumber $ws.Range("V6") " -122.419791"

# Row 7 - man3.jpg observation, lat/lon
$ws.Range("T7").Value = "observations/portrait/man3.jpg"
$ws.Range("T7").Font.Color = 0
Set-TextLooksLikeNumber $ws.Range("U7") " 37.767076"
Set-TextLooksLikeNumber $ws.Range("V7") " -122.419918"

# Row 8 - man4.jpg observation, lat/lon
$ws.Range("T8").Value = "observations/portrait/man4.jpg"
$ws.Range("T8").Font.Color = 0
Set-TextLooksLikeNumber $ws.Range("U8") " 37.756364"
Set-TextLooksLikeNumber $ws.Range("V8") " -122.421321"

# Row 9 - woman2.jpg observation, lat/lon (numeric)
$ws.Range("T9").Value = "observations/portrait/woman2.jpg"
$ws.Range("T9").Font.Color = 0
$ws.Range("U9").Value = 37.785252
$ws.Range("V9").Value = -122.403587

# Row 10 - woman3.jpg observation, lat/lon (numeric)
$ws.Range("T10").Value = "observations/portrait/woman3.jpg"
$ws.Range("T10").Font.Color = 0
$ws.Range("U10").Value = 37.785800000000002
$ws.Range("V10").Value = -122.404113
